# Weekly update: shift existing records down by inserting two new rows
# (a fresh "Primera"/"Segunda" pair) at the top of the data block (row 18),
# pushing the older rows down two positions. The two rows that fall off the
# bottom of the original block reappear at rows 41-42 automatically as part
# of the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 18; everything from row 18 down moves to row 20+.
$ws.Rows.Item(18).EntireRow.Insert()
$ws.Rows.Item(18).EntireRow.Insert()

# New row 18: "Primera" quality record for the latest week.
$ws.Cells.Item(18, 1).Value = 11
$ws.Cells.Item(18, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(18, 3).Value = "Bíobío"
$ws.Cells.Item(18, 4).Value = [DateTime]"2022-04-14"
$ws.Cells.Item(18, 5).Value = 8
$ws.Cells.Item(18, 6).Value = 100112043
$ws.Cells.Item(18, 7).Value = "Pepino dulce"
$ws.Cells.Item(18, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 200
$ws.Cells.Item(18, 11).Value = 15000
$ws.Cells.Item(18, 12).Value = 16000
$ws.Cells.Item(18, 13).Value = 15500
$ws.Cells.Item(18, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 861
$ws.Cells.Item(18, 17).Value = 18
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# New row 19: "Segunda" quality record for the same latest week.
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(19, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(19, 3).Value = "Bíobío"
$ws.Cells.Item(19, 4).Value = [DateTime]"2022-04-14"
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(19, 6).Value = 100112043
$ws.Cells.Item(19, 7).Value = "Pepino dulce"
$ws.Cells.Item(19, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(19, 9).Value = "Segunda"
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 13000
$ws.Cells.Item(19, 12).Value = 13000
$ws.Cells.Item(19, 13).Value = 13000
$ws.Cells.Item(19, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(19, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(19, 16).Value = 722
$ws.Cells.Item(19, 17).Value = 18
$ws.Cells.Item(19, 18).Value = "Hortaliza"
